# Update "想去人数" (interested-count) figures on the "展览" and "全部类型"
# sheets to reflect newly scraped totals.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 1770
    $ws.Range("F3").Value = 8072
    $ws.Range("F4").Value = 187
    $ws.Range("F5").Value = 289
}
